$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ticket synced from Mobile App (new date, new ticket_id, job reassigned) ---
$ws.Range("A2").Value = (Get-Date -Year 2025 -Month 12 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B2").Value = 413345
$ws.Range("C2").Value = "normal"
$ws.Range("D2").Value = "Instalación"
$ws.Range("E2").Value = "Antena GPS"
$ws.Range("F2").Value = "AMERICO VESPUCIO 2341 PUDAHUEL"
$ws.Range("G2").Value = "PUDAHUEL"
$ws.Range("H2").Value = "Región Metropolitana de Santiago."
$ws.Range("I2").Value = "Juan Perez"
$ws.Range("J2").Value = "HCCR34"
$ws.Range("K2").Value = "YALA"

# --- Row 3: ticket synced from Mobile App (new date, new ticket_id, job reassigned) ---
$ws.Range("A3").Value = (Get-Date -Year 2025 -Month 12 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B3").Value = 413235
$ws.Range("C3").Value = "normal"
$ws.Range("D3").Value = "Instalación"
$ws.Range("E3").Value = "Antena GPS"
$ws.Range("F3").Value = "AMERICO VESPUCIO 2341 PUDAHUEL"
$ws.Range("G3").Value = "PUDAHUEL"
$ws.Range("H3").Value = "Región Metropolitana de Santiago."
$ws.Range("I3").Value = "Pedro Pascal"
$ws.Range("J3").Value = "BSBJ92"
$ws.Range("K3").Value = "VALU"

# --- Backend plan reset: wipe out the two now-stale planning rows ---
$ws.Range("A4:K5").ClearContents() | Out-Null

# --- Extend the reset plan with two fresh blank rows, keeping column A's
#     date-formatted style consistent with the rest of the plan ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths: best-fit to the new (wider) plan contents ---
$ws.Range("A1").ColumnWidth = 9.083334333333331
$ws.Range("B1").ColumnWidth = 6.916667666666667
$ws.Range("C1").ColumnWidth = 8.083334333333331
$ws.Range("D1").ColumnWidth = 10.083334333333331
$ws.Range("E1").ColumnWidth = 20.250001
$ws.Range("F1").ColumnWidth = 30.250001
$ws.Range("G1").ColumnWidth = 8.750001
$ws.Range("H1").ColumnWidth = 28.250001
$ws.Range("I1").ColumnWidth = 13.250001
$ws.Range("J1").ColumnWidth = 22.750001
$ws.Range("K1").ColumnWidth = 15.916667666666665

# --- View state: keep selection on K4 (matches the saved selection) ---
$ws.Range("K4").Select() | Out-Null
